$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be treated as text so values such as
# "220.22", "1.004" or "26.430.09" are not re-interpreted as numbers/dates
# (the source data always stores these as plain text).
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "26.430.09"
$ws.Range("E2").Value = "  +1.68%  "
$ws.Range("D3").Value = "1.669.77"
$ws.Range("E3").Value = "  +1.35%  "
$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  +0.39%  "
$ws.Range("D5").Value = "220.22"
$ws.Range("E5").Value = "  +2.94%  "
$ws.Range("D6").Value = "0.5260"
$ws.Range("E6").Value = "  +0.69%  "
$ws.Range("E7").Value = "  +0.44%  "
$ws.Range("D8").Value = "0.2670"
$ws.Range("E8").Value = "  +2.26%  "
$ws.Range("D9").Value = "0.06363"
$ws.Range("E9").Value = "  +0.25%  "
$ws.Range("D10").Value = "21.71"
$ws.Range("E10").Value = "  +4.54%  "
$ws.Range("D11").Value = "0.07786"
$ws.Range("E11").Value = "  +1.37%  "
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").Value = "4.466"
$ws.Range("E12").Value = "  +0.87%  "
$ws.Range("B13").Value = "WrappedEther"
$ws.Range("C13").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D13").Value = "1.653.48"
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "0.5526"
$ws.Range("E14").Value = "  +0.37%  "
$ws.Range("D15").Value = "0.0₅8263"
$ws.Range("D16").Value = "65.54"
$ws.Range("E16").Value = "  +1.31%  "
$ws.Range("D17").Value = "26.439.97"
$ws.Range("E17").Value = "  +1.77%  "
$ws.Range("E18").Value = "  +0.34%  "
$ws.Range("D19").Value = "4.738"
$ws.Range("E19").Value = "  +0.43%  "
$ws.Range("D20").Value = "193.70"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("D21").Value = "10.26"
$ws.Range("E21").Value = "  +0.59%  "
$ws.Range("D22").Value = "6.272"
$ws.Range("D23").Value = "1.005"
$ws.Range("E23").Value = "  +0.49%  "
$ws.Range("B24").Value = "Stellar"
$ws.Range("C24").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D24").Value = "0.1262"
$ws.Range("E24").Value = "  +1.28%  "
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D25").Value = "138.85"
$ws.Range("E25").Value = "  -3.04%  "
$ws.Range("D26").Value = "7.381"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "16.21"
$ws.Range("E27").Value = "  +1.44%  "
$ws.Range("D28").Value = "1.420"
$ws.Range("E28").Value = "  +0.96%  "
$ws.Range("D29").Value = "0.06127"
$ws.Range("E29").Value = "  +3.26%  "
$ws.Range("E30").Value = "  +2.64%  "
$ws.Range("D31").Value = "3.589"
$ws.Range("E31").Value = "  +4.83%  "
$ws.Range("D32").Value = "3.397"
$ws.Range("E32").Value = "  -0.30%  "
$ws.Range("D33").Value = "1.679"
$ws.Range("E33").Value = "  +2.24%  "
$ws.Range("D34").Value = "1.003"
$ws.Range("E34").Value = "  +1.06%  "
$ws.Range("D35").Value = "2.423"
$ws.Range("E35").Value = "  +1.34%  "
$ws.Range("D36").Value = "0.6052"
$ws.Range("E36").Value = "  +7.42%  "
$ws.Range("D37").Value = "2.766"
$ws.Range("E37").Value = "  +0.83%  "
$ws.Range("D38").Value = "0.01612"
$ws.Range("E38").Value = "  +0.46%  "
$ws.Range("D39").Value = "6.033"
$ws.Range("E39").Value = "  +2.91%  "
$ws.Range("D40").Value = "1.088.98"
$ws.Range("E40").Value = "  +6.19%  "
$ws.Range("D41").Value = "0.8590"
$ws.Range("E41").Value = "  +0.45%  "
$ws.Range("E42").Value = "  +0.27%  "
$ws.Range("D43").Value = "100.60"
$ws.Range("E43").Value = "  +1.70%  "
$ws.Range("D44").Value = "1.811.52"
$ws.Range("E44").Value = "  +0.98%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "57.88"
$ws.Range("E45").Value = "  +4.03%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₈108"
$ws.Range("E46").Value = "  +0.45%  "
$ws.Range("B47").Value = "Frax"
$ws.Range("C47").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D47").Value = "1.004"
$ws.Range("E47").Value = "  +0.31%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "8.123"
$ws.Range("E48").Value = "  +0.85%  "
$ws.Range("D49").Value = "0.05202"
$ws.Range("E49").Value = "  +1.20%  "
$ws.Range("D50").Value = "1.481"
$ws.Range("E50").Value = "  +7.65%  "
$ws.Range("D51").Value = "0.4229"
$ws.Range("E51").Value = "  +0.51%  "

# Restore the default ("Normal") cell style on the Price column so no stray
# formatting is left behind now that the text values have been written.
$ws.Range("D2:D51").Style = "Normal"
